$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff inserts one new data row right before the current row 112 (Berenjena /
# Macroferia Regional de Talca sheet), which pushes the existing rows 112-144 down
# to 113-145 (dimension grows from A1:R144 to A1:R145) and fills the freshly
# opened row 112 with a brand-new record.

$ws.Rows("112:112").Insert()

$ws.Cells.Item(112, 1).Value  = 5
$ws.Cells.Item(112, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value  = "Maule"
$ws.Range("D112").Value       = 44985
$ws.Cells.Item(112, 5).Value  = 7
$ws.Cells.Item(112, 6).Value  = 100112001
$ws.Cells.Item(112, 7).Value  = "Berenjena"
$ws.Cells.Item(112, 8).Value  = "Sin especificar"
$ws.Cells.Item(112, 9).Value  = "Primera"
$ws.Cells.Item(112, 10).Value = 200
$ws.Cells.Item(112, 11).Value = 8000
$ws.Cells.Item(112, 12).Value = 8000
$ws.Cells.Item(112, 13).Value = 8000
$ws.Cells.Item(112, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(112, 15).Value = "Región del Maule"
$ws.Cells.Item(112, 16).Value = 160
$ws.Cells.Item(112, 17).Value = 50
$ws.Cells.Item(112, 18).Value = "Hortaliza"
